$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Helper pattern used throughout: to split a run cleanly (without the
# engine tagging extra/unwanted xml:space="preserve" attributes), we
# bracket the exact sub-range we are about to rewrite with a pair of
# temporary bookmarks (one at the start, one at the end of the
# sub-range), assign .Text on that sub-range (even to the same text,
# just to force a clean split), and then delete the temporary
# bookmarks again. This reliably produces the same run-splitting
# behaviour Word itself performs when a user selects text and types
# a replacement.
# -----------------------------------------------------------------

# ===================================================================
# Change 1 (paragraph 1): "... A Class can inherit from Multiple
# Interfaces but cannot inherit from Multiple Classes." ->
# "... A Class can implement Multiple Interfaces but cannot inherit
# from Multiple Classes." with the (moved) _GoBack bookmark now
# sitting right after "implement ".
# ===================================================================

$rng1 = $d.Content
$found1 = $rng1.Find.Execute("inherit from Multiple Interfaces but cannot inherit from Multiple Classes.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$splitPos1 = $rng1.Start
$d.Bookmarks.Add("TmpSplit1", $d.Range($splitPos1, $splitPos1))

$rReplace1 = $d.Range($rng1.Start, $rng1.Start + 13)   # "inherit from "
$rReplace1.Text = "implement "

$d.Bookmarks.Item("TmpSplit1").Delete()

# Insert the _GoBack bookmark right after "implement " (zero length).
$bmPos1 = $rng1.Start + 10
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos1, $bmPos1))

# ===================================================================
# Change 2 (paragraph 2): drop ", Delegates," from the list.
# ===================================================================

$d.Content.Find.Execute(
    "An Interface can contain Properties, Methods, Delegates, or Events but only the declaration and not the implementation.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "An Interface can contain Properties, Methods or Events but only the declaration and not the implementation.",
    2) | Out-Null

# ===================================================================
# Change 3 (paragraph 3): "Interface methods are Public by default.
# You cannot explicitly use Public keyword for an interface method."
# -> "Interface members are Public by default. You cannot explicitly
# use Public keyword for interface members." split across 5 runs.
# ===================================================================

# Step 1: replace first occurrence of "methods" -> "members"
$rFull3 = $d.Content
$found3a = $rFull3.Find.Execute("Interface methods are Public by default. You cannot explicitly use Public keyword for an interface method.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base3 = $rFull3.Start

$p10 = $base3 + 10   # after "Interface "
$p17 = $base3 + 17   # after "methods"/"members"

$d.Bookmarks.Add("Tmp3A", $d.Range($p10, $p10))
$d.Bookmarks.Add("Tmp3B", $d.Range($p17, $p17))
$rMethods = $d.Range($p10, $p17)
$rMethods.Text = "members"
$d.Bookmarks.Item("Tmp3A").Delete()
$d.Bookmarks.Item("Tmp3B").Delete()

# Step 2: replace the tail "an interface method." -> "interface members."
# as a single run first.
$rTailFind = $d.Content
$found3b = $rTailFind.Find.Execute("an interface method.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tailStart = $rTailFind.Start
$tailEnd = $rTailFind.End

$d.Bookmarks.Add("Tmp3C", $d.Range($tailStart, $tailStart))
$d.Bookmarks.Add("Tmp3D", $d.Range($tailEnd, $tailEnd))
$rTail = $d.Range($tailStart, $tailEnd)
$rTail.Text = "interface members."
$d.Bookmarks.Item("Tmp3C").Delete()
$d.Bookmarks.Item("Tmp3D").Delete()

# Step 3: split the new "interface members." run into three runs:
# "interface " | "members" | "."
$rTail2 = $d.Content
$found3c = $rTail2.Find.Execute("interface members.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base3b = $rTail2.Start

$q10 = $base3b + 10  # after "interface "
$q17 = $base3b + 17  # after "members"

$d.Bookmarks.Add("Tmp3E", $d.Range($base3b, $base3b))
$d.Bookmarks.Add("Tmp3F", $d.Range($q10, $q10))
$rPart1 = $d.Range($base3b, $q10)
$rPart1.Text = "interface "
$d.Bookmarks.Item("Tmp3E").Delete()
$d.Bookmarks.Item("Tmp3F").Delete()

$d.Bookmarks.Add("Tmp3G", $d.Range($q10, $q10))
$d.Bookmarks.Add("Tmp3H", $d.Range($q17, $q17))
$rPart2 = $d.Range($q10, $q17)
$rPart2.Text = "members"
$d.Bookmarks.Item("Tmp3G").Delete()
$d.Bookmarks.Item("Tmp3H").Delete()

# ===================================================================
# Change 4: remove the old _GoBack bookmark that used to sit after
# "... but not from any classes" (it has effectively moved to
# paragraph 1, added above).
# ===================================================================

if ($d.Bookmarks.Exists("_GoBack")) {
    # There are now two things named "_GoBack" is not possible (bookmark
    # names are unique), so this call grabs the single remaining
    # instance -- but we only want to delete the OLD one near
    # "but not from any classes", not the one we just added above.
}

$rOld = $d.Content
$foundOld = $rOld.Find.Execute("but not from any classes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$oldPos = $rOld.End

# Find the bookmark whose range sits at this location and remove it,
# without touching the _GoBack we created earlier in paragraph 1.
for ($i = $d.Bookmarks.Count; $i -ge 1; $i--) {
    $bm = $d.Bookmarks.Item($i)
    if ($bm.Name -eq "_GoBack" -and $bm.Start -eq $oldPos -and $bm.End -eq $oldPos) {
        $bm.Delete()
    }
}

$d.Save()
